# Generate Report for Handback
# Updates the localization-status workbook so the zh-cn and de-de handback
# files are reported as freshly handed back / in sync, clearing the stale
# "not the latest" error and refreshing the handback timestamps, and widens
# the Status / Error Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# --- zh-cn sheet: Status, Latest Handback DateTime, Error Detail ---
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-08-21 21:00:19"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: Status, Latest Handback DateTime, Error Detail ---
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-08-21 21:00:26"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments to accommodate the longer status text and
#     the now-empty error column (nearest values on Excel's column-width
#     grid to the authored target widths) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
